$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before row 330 (existing rows 330:345 shift down to 334:349).
$ws.Range("A330:A333").EntireRow.Insert()

# Row 330 - Especial, M=1000
$ws.Cells.Item(330,1).Value2  = 5
$ws.Cells.Item(330,2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(330,3).Value2  = "Maule"
$ws.Cells.Item(330,4).Value2  = 44509
$ws.Cells.Item(330,5).Value2  = 7
$ws.Cells.Item(330,6).Value2  = "Fruta"
$ws.Cells.Item(330,7).Value2  = 100101
$ws.Cells.Item(330,8).Value2  = "Berries"
$ws.Cells.Item(330,9).Value2  = 100112025
$ws.Cells.Item(330,10).Value2 = "Frutilla"
$ws.Cells.Item(330,11).Value2 = "Sin especificar"
$ws.Cells.Item(330,12).Value2 = "Especial"
$ws.Cells.Item(330,13).Value2 = 1000
$ws.Cells.Item(330,14).Value2 = 8000
$ws.Cells.Item(330,15).Value2 = 8000
$ws.Cells.Item(330,16).Value2 = 8000
$ws.Cells.Item(330,17).Value2 = "$/bandeja 7 kilos"
$ws.Cells.Item(330,18).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(330,19).Value2 = 1143
$ws.Cells.Item(330,20).Value2 = 7

# Row 331 - Especial, M=100
$ws.Cells.Item(331,1).Value2  = 5
$ws.Cells.Item(331,2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(331,3).Value2  = "Maule"
$ws.Cells.Item(331,4).Value2  = 44509
$ws.Cells.Item(331,5).Value2  = 7
$ws.Cells.Item(331,6).Value2  = "Fruta"
$ws.Cells.Item(331,7).Value2  = 100101
$ws.Cells.Item(331,8).Value2  = "Berries"
$ws.Cells.Item(331,9).Value2  = 100112025
$ws.Cells.Item(331,10).Value2 = "Frutilla"
$ws.Cells.Item(331,11).Value2 = "Sin especificar"
$ws.Cells.Item(331,12).Value2 = "Especial"
$ws.Cells.Item(331,13).Value2 = 100
$ws.Cells.Item(331,14).Value2 = 9000
$ws.Cells.Item(331,15).Value2 = 9000
$ws.Cells.Item(331,16).Value2 = 9000
$ws.Cells.Item(331,17).Value2 = "$/caja 7 kilos"
$ws.Cells.Item(331,18).Value2 = "Región del Maule"
$ws.Cells.Item(331,19).Value2 = 1286
$ws.Cells.Item(331,20).Value2 = 7

# Row 332 - Primera, M=50
$ws.Cells.Item(332,1).Value2  = 5
$ws.Cells.Item(332,2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(332,3).Value2  = "Maule"
$ws.Cells.Item(332,4).Value2  = 44509
$ws.Cells.Item(332,5).Value2  = 7
$ws.Cells.Item(332,6).Value2  = "Fruta"
$ws.Cells.Item(332,7).Value2  = 100101
$ws.Cells.Item(332,8).Value2  = "Berries"
$ws.Cells.Item(332,9).Value2  = 100112025
$ws.Cells.Item(332,10).Value2 = "Frutilla"
$ws.Cells.Item(332,11).Value2 = "Sin especificar"
$ws.Cells.Item(332,12).Value2 = "Primera"
$ws.Cells.Item(332,13).Value2 = 50
$ws.Cells.Item(332,14).Value2 = 8000
$ws.Cells.Item(332,15).Value2 = 8000
$ws.Cells.Item(332,16).Value2 = 8000
$ws.Cells.Item(332,17).Value2 = "$/caja 7 kilos"
$ws.Cells.Item(332,18).Value2 = "Región del Maule"
$ws.Cells.Item(332,19).Value2 = 1143
$ws.Cells.Item(332,20).Value2 = 7

# Row 333 - Segunda, M=200
$ws.Cells.Item(333,1).Value2  = 5
$ws.Cells.Item(333,2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(333,3).Value2  = "Maule"
$ws.Cells.Item(333,4).Value2  = 44509
$ws.Cells.Item(333,5).Value2  = 7
$ws.Cells.Item(333,6).Value2  = "Fruta"
$ws.Cells.Item(333,7).Value2  = 100101
$ws.Cells.Item(333,8).Value2  = "Berries"
$ws.Cells.Item(333,9).Value2  = 100112025
$ws.Cells.Item(333,10).Value2 = "Frutilla"
$ws.Cells.Item(333,11).Value2 = "Sin especificar"
$ws.Cells.Item(333,12).Value2 = "Segunda"
$ws.Cells.Item(333,13).Value2 = 200
$ws.Cells.Item(333,14).Value2 = 5000
$ws.Cells.Item(333,15).Value2 = 5000
$ws.Cells.Item(333,16).Value2 = 5000
$ws.Cells.Item(333,17).Value2 = "$/bandeja 7 kilos"
$ws.Cells.Item(333,18).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(333,19).Value2 = 714
$ws.Cells.Item(333,20).Value2 = 7
